$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Priority changed from "low" to "ht" for the four "Ready for handoff" rows
# (rows 4-7) in both the zh-cn and de-de localization-status sheets.
$wsZh.Range("E4:E7").Value = "ht"
$wsDe.Range("E4:E7").Value = "ht"

# Latest Handoff Datetime refreshed for the same rows, per language.
$wsZh.Range("H4:H7").Value = "2016-09-02 00:40:38"
$wsDe.Range("H4:H7").Value = "2016-09-02 00:40:43"

# The Overview sheet's "Latest HO Xliff Generate Date" mirrors the newest
# per-language handoff datetime (de-de, which is later than zh-cn here).
$wsOverview.Range("G4:G7").Value = "2016-09-02 00:40:43"
